# Applies the "Batterywise analysis" relabeling/reshuffling edit to the
# "Analysis Results" sheet.
#
# The edit: a batch of row labels in column A get clarified with units
# (e.g. "Peak Power" -> "Peak Power(kW)"), a couple of rows get swapped
# (Starting/Ending SoC %, Lowest/Highest Cell Voltage, lowest/highest cell
# temp), several numeric values in column B are recomputed, and two new
# rows (42 "Time spent in 70-80 km/h" and 43 "Time spent in 80-90 km/h")
# are appended at the bottom, growing the used range from A1:B42 to A1:B43.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6-7: Starting/Ending SoC (%) values swap (labels stay the same).
$ws.Cells.Item(6, 2).Value = 99
$ws.Cells.Item(7, 2).Value = 18

# Row 8-10: clarify labels with units.
$ws.Cells.Item(8, 1).Value = "Total distance covered (km)"
$ws.Cells.Item(9, 1).Value = "Total energy consumption(WH/KM)"
$ws.Cells.Item(10, 1).Value = "Total SOC consumed(%)"

# Row 12-14: clarify labels with units.
$ws.Cells.Item(12, 1).Value = "Peak Power(kW)"
$ws.Cells.Item(13, 1).Value = "Average Power(kW)"
$ws.Cells.Item(14, 1).Value = "Total Energy Regenerated(kWh)"

# Row 15: clarify label, flip sign of value.
$ws.Cells.Item(15, 1).Value = "Regenerative Effectiveness(%)"
$ws.Cells.Item(15, 2).Value = 0.004523068710953785

# Row 16-17: Lowest/Highest Cell Voltage swap (values move with new labels).
$ws.Cells.Item(16, 1).Value = "Highest Cell Voltage(V)"
$ws.Cells.Item(16, 2).Value = 3.394
$ws.Cells.Item(17, 1).Value = "Lowest Cell Voltage(V)"
$ws.Cells.Item(17, 2).Value = 3.094

# Row 18-20: clarify labels with units.
$ws.Cells.Item(18, 1).Value = "Difference in Cell Voltage(V)"
$ws.Cells.Item(19, 1).Value = "Minimum Temperature(C)"
$ws.Cells.Item(20, 1).Value = "Maximum Temperature(C)"

# Row 21: clarify label, fill in previously-empty value.
$ws.Cells.Item(21, 1).Value = "Difference in Temperature(C)"
$ws.Cells.Item(21, 2).Value = 9

# Row 22-27: clarify labels with units.
$ws.Cells.Item(22, 1).Value = "Maximum Fet Temperature-BMS(C)"
$ws.Cells.Item(23, 1).Value = "Maximum Afe Temperature-BMS(C)"
$ws.Cells.Item(24, 1).Value = "Maximum PCB Temperature-BMS(C)"
$ws.Cells.Item(25, 1).Value = "Maximum MCU Temperature(C)"
$ws.Cells.Item(26, 1).Value = "Maximum Motor Temperature(C)"
$ws.Cells.Item(27, 1).Value = "Abnormal Motor Temperature Detected(C)"

# Row 28-29: lowest/highest cell temp swap (values move with new labels).
$ws.Cells.Item(28, 1).Value = "highest cell temp(C)"
$ws.Cells.Item(29, 1).Value = "lowest cell temp(C)"

# Row 30: clarify label with units.
$ws.Cells.Item(30, 1).Value = "Difference between Highest and Lowest Cell Temperature at 100% SOC(C)"

# Row 31: relabeled to Battery Voltage(V), new value.
$ws.Cells.Item(31, 1).Value = "Battery Voltage(V)"
$ws.Cells.Item(31, 2).Value = 53

# Row 32: relabeled to Total energy charged(kWh), new value.
$ws.Cells.Item(32, 1).Value = "Total energy charged(kWh)"
$ws.Cells.Item(32, 2).Value = 1.751019874166667

# Row 33: relabeled to Electricity consumption units(kW), new value.
$ws.Cells.Item(33, 1).Value = "Electricity consumption units(kW)"
$ws.Cells.Item(33, 2).Value = 0.0000001382980976658347

# Row 34: relabeled to Idling time percentage, new value.
$ws.Cells.Item(34, 1).Value = "Idling time percentage"
$ws.Cells.Item(34, 2).Value = 10.21450141651879

# Row 35: relabeled to Time spent in 0-10 km/h, new value.
$ws.Cells.Item(35, 1).Value = "Time spent in 0-10 km/h"
$ws.Cells.Item(35, 2).Value = 4.321160611438

# Row 36: relabeled to Time spent in 10-20 km/h, value unchanged.
$ws.Cells.Item(36, 1).Value = "Time spent in 10-20 km/h"
$ws.Cells.Item(36, 2).Value = 2.375393045048411

# Row 37: relabeled to Time spent in 20-30 km/h, value unchanged.
$ws.Cells.Item(37, 1).Value = "Time spent in 20-30 km/h"
$ws.Cells.Item(37, 2).Value = 5.220883534136546

# Row 38: relabeled to Time spent in 30-40 km/h, new value.
$ws.Cells.Item(38, 1).Value = "Time spent in 30-40 km/h"
$ws.Cells.Item(38, 2).Value = 10.88073223125058

# Row 39: relabeled to Time spent in 40-50 km/h, new value.
$ws.Cells.Item(39, 1).Value = "Time spent in 40-50 km/h"
$ws.Cells.Item(39, 2).Value = 11.01771426792441

# Row 40: relabeled to Time spent in 50-60 km/h, new value.
$ws.Cells.Item(40, 1).Value = "Time spent in 50-60 km/h"
$ws.Cells.Item(40, 2).Value = 51.5457177547399

# Row 41: relabeled to Time spent in 60-70 km/h, new value.
$ws.Cells.Item(41, 1).Value = "Time spent in 60-70 km/h"
$ws.Cells.Item(41, 2).Value = 4.342953208181564

# Row 42: relabeled to Time spent in 70-80 km/h, new value.
$ws.Cells.Item(42, 1).Value = "Time spent in 70-80 km/h"
$ws.Cells.Item(42, 2).Value = 0

# Row 43 (new): Time spent in 80-90 km/h.
$ws.Cells.Item(43, 1).Value = "Time spent in 80-90 km/h"
$ws.Cells.Item(43, 2).Value = 0
